# "Nueva tarea en el productBacklog"
# Inserts a new Product Backlog Item row ("Ejecutar un ejemplo de Prueba e
# instalar en apk") as item #2, pushing the previously existing items 2-6
# down to items 3-7, and tidies up a couple of rows at the bottom of the
# sheet accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")
$ws.Activate()

$xlPasteAll = -4104
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Shift the existing task rows (4-8) down to (5-9), bottom-up so we
#    never overwrite data before it has been copied. This carries both
#    the values and the existing formatting along with it.
# ---------------------------------------------------------------------
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial($xlPasteAll)

$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial($xlPasteAll)

$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial($xlPasteAll)

$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial($xlPasteAll)

$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial($xlPasteAll)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Row 9 used to be a blank spacer row (now holding what was row 8's
#    data thanks to the shift above) - its C:E cells should look like
#    plain body cells instead of carrying the old spacer shading.
# ---------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("C9:E9").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Turn row 4 into the brand-new backlog item. Clear it completely
#    first (values + formatting) and retype it, which is what the
#    author actually did for the "Tarea" cell (it ends up with no
#    explicit style at all, i.e. the default style).
# ---------------------------------------------------------------------
$ws.Range("A4:E4").Clear()

$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").Value = 2

$ws.Range("B4").Value = "Ejecutar un ejemplo de Prueba e instalar en apk"

$ws.Range("C5:E5").Copy()
$ws.Range("C4:E4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. A few knock-on formatting touch-ups further down the sheet.
# ---------------------------------------------------------------------
# Row 10 (first blank row after the items) - A column takes on the look
# previously used by the item rows' first column.
$ws.Range("A5").Copy()
$ws.Range("A10").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Row 11 - B column takes on the item-row shading as well.
$ws.Range("B5").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Row 12 - B column reverts back to the plain body-cell look.
$ws.Range("C3").Copy()
$ws.Range("B12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Reshuffle the tail of the sheet: what used to be a single all-blank
#    row 20 becomes a short row 20 (only column A, keeping the boxed
#    look) followed by a new, fully blank row 21. Row 19 (B:E) loses its
#    shading to match the new plain trailing rows.
# ---------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B19:E19").ClearContents()
$ws.Range("A20").Copy()
$ws.Range("B19:E19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B20:E20").Clear()
$ws.Range("A11").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 6. Restore the selection to match where the author ended up editing.
# ---------------------------------------------------------------------
$ws.Range("B9").Select()
